$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.522.27"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").Value = "1.576.02"
$ws.Range("E3").Value = "  -3.46%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.76%  "

$ws.Range("E6").Value = "  -3.36%  "

$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.251"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0588"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.94%  "

$ws.Range("E11").Value = "  -2.19%  "

$ws.Range("D12").Value = "1.798.01"
$ws.Range("E12").Value = "  -3.57%  "

$ws.Range("D13").Value = "1.565.89"
$ws.Range("E13").Value = "  -4.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.93%  "

$ws.Range("E15").Value = "  -6.55%  "

$ws.Range("D16").Value = "27.500.20"
$ws.Range("E16").Value = "  -1.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.62%  "

$ws.Range("D19").Value = "0.0₃0689"
$ws.Range("E19").Value = "  -4.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.99%  "

$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.83%  "

$ws.Range("E23").Value = "  -6.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.79%  "

$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("E27").Value = "  -3.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.46%  "

$ws.Range("E29").Value = "  -4.80%  "

$ws.Range("E30").Value = "  -2.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0462"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.82%  "

$ws.Range("E32").Value = "  -5.59%  "

$ws.Range("D33").Value = "1.357.89"
$ws.Range("E33").Value = "  -2.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.93%  "

$ws.Range("E35").Value = "  -5.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.964"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.55%  "

$ws.Range("E37").Value = "  -1.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0164"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.534"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.72%  "

$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.970"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.40%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.22%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.27%  "

$ws.Range("D47").Value = "1.710.61"
$ws.Range("E47").Value = "  -3.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.99%  "

$ws.Range("D49").Value = "0.0₆01000"
$ws.Range("E49").Value = "  -3.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0967"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.83%  "

$ws.Range("E51").Value = "  -1.84%  "

